$d = $word.ActiveDocument

# Step 1: "Kierownik studiów pierwszego stopnia" -> "Prodziekan ds. studenckich "
$null = $d.Content.Find.Execute("Kierownik studiów pierwszego stopnia", $false, $false, $false, $false, $false, $true, 1, $false, "Prodziekan ds. studenckich ", 2)

# Step 2: restructure the two paragraphs that followed ("na kierunku Informatyka" and
# "dr " + "Jakub Zygadło" + bookmark) into three paragraphs:
#   - a blank-indent line with a long run of spaces + "Wydziału Matematyki i Informatyki "
#   - "dr hab. Piotr Niemiec"
#   - an (otherwise empty) paragraph carrying only the _GoBack bookmark
$pA = $d.Paragraphs.Item(15)
$pB = $d.Paragraphs.Item(16)
$combined = $d.Range($pA.Range.Start, $pB.Range.End)
$combined.InsertXML("<w:p xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main""><w:pPr><w:pStyle w:val=""Default""/><w:rPr><w:rFonts w:asciiTheme=""minorHAnsi"" w:hAnsiTheme=""minorHAnsi""/><w:sz w:val=""16""/><w:szCs w:val=""16""/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme=""minorHAnsi"" w:hAnsiTheme=""minorHAnsi""/><w:i/><w:iCs/><w:sz w:val=""23""/><w:szCs w:val=""23""/></w:rPr><w:t xml:space=""preserve"">                                                                                               Wydziału Matematyki i Informatyki </w:t></w:r></w:p><w:p xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main""><w:pPr><w:pStyle w:val=""Default""/><w:ind w:left=""4248"" w:firstLine=""708""/><w:rPr><w:rFonts w:asciiTheme=""minorHAnsi"" w:hAnsiTheme=""minorHAnsi""/><w:i/><w:iCs/><w:sz w:val=""23""/><w:szCs w:val=""23""/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme=""minorHAnsi"" w:hAnsiTheme=""minorHAnsi""/><w:i/><w:iCs/><w:sz w:val=""23""/><w:szCs w:val=""23""/></w:rPr><w:t>dr hab. Piotr Niemiec</w:t></w:r></w:p><w:p xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main""><w:pPr><w:pStyle w:val=""Default""/><w:ind w:left=""4248"" w:firstLine=""708""/><w:rPr><w:rFonts w:asciiTheme=""minorHAnsi"" w:hAnsiTheme=""minorHAnsi""/><w:i/><w:iCs/><w:sz w:val=""23""/><w:szCs w:val=""23""/></w:rPr></w:pPr><w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/></w:p>")

Write-Host "done"
